# Prepend a label ("Ano "/"Intervalo ") to the header row (row 1) cells
# so that, once imported into Power BI, the first row can automatically
# be promoted to become the table header.

$wb = $excel.ActiveWorkbook

# Sheets 1, 2, 3 and 5 use year headers (2015 / 2030 / 2040 / 2050) -> "Ano <year>"
$anoSheets = @(1, 2, 3, 5)
foreach ($idx in $anoSheets) {
    $ws = $wb.Worksheets.Item($idx)
    foreach ($col in @("B", "C", "D", "E")) {
        $cell = $ws.Range($col + "1")
        $cell.Value = "Ano " + $cell.Value2
    }
}

# Sheet 4 uses period headers (2015 / 2015-2030 / 2031-2040 / 2041-2050) -> "Intervalo <period>"
$ws4 = $wb.Worksheets.Item(4)
foreach ($col in @("B", "C", "D", "E")) {
    $cell = $ws4.Range($col + "1")
    $cell.Value = "Intervalo " + $cell.Value2
}

# Sheet 6 only has a single year column (B1) -> "Ano <year>"
$ws6 = $wb.Worksheets.Item(6)
$cellB1 = $ws6.Range("B1")
$cellB1.Value = "Ano " + $cellB1.Value2
